$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 337.301872878912
$ws.Range("F2").Value = 175.0122316609742
$ws.Range("G2").Value = -0.0000000000858635347429663
$ws.Range("K2").Value = 0.4525505391315189
$ws.Range("C3").Value = 1177.61873778985
$ws.Range("F3").Value = 397.1059129019827
$ws.Range("G3").Value = 245.6123408855054
$ws.Range("K3").Value = 1.92747852380969
$ws.Range("C4").Value = 2271.270931647058
$ws.Range("F4").Value = 648.3402844044327
$ws.Range("G4").Value = 802.5533206762635
$ws.Range("K4").Value = 4.617308479895989
$ws.Range("C5").Value = 3334.050946095882
$ws.Range("F5").Value = 897.763784803767
$ws.Range("G5").Value = 1714.343376201689
$ws.Range("K5").Value = 8.589657843470457
$ws.Range("C6").Value = 4316.75863830004
$ws.Range("F6").Value = 1156.953952186595
$ws.Range("K6").Value = 13.86047259264091
$ws.Range("C7").Value = 5229.865727713754
$ws.Range("F7").Value = 1425.831136761117
$ws.Range("G7").Value = 4607.225239504814
$ws.Range("K7").Value = 20.47561536281522
$ws.Range("C8").Value = 6072.574524034305
$ws.Range("F8").Value = 1701.377148470918
$ws.Range("G8").Value = 6614.683984933854
$ws.Range("K8").Value = 28.47542750499196
$ws.Range("C9").Value = 6765.849125159718
$ws.Range("F9").Value = 1987.147194107323
$ws.Range("G9").Value = 9007.839350411174
$ws.Range("K9").Value = 37.79535684183598
$ws.Range("C10").Value = 7050.874324237787
$ws.Range("F10").Value = 2297.737444150939
$ws.Range("G10").Value = 11649.68579295754
$ws.Range("K10").Value = 48.7581256366926
$ws.Range("C11").Value = 7838.960082260692
$ws.Range("F11").Value = 2614.389296533142
$ws.Range("G11").Value = 14983.8166790652
$ws.Range("K11").Value = 61.08817982711074
$ws.Range("C12").Value = 7842.34378029345
$ws.Range("F12").Value = -2614.439928077682
$ws.Range("G12").Value = 18646.69199200964
$ws.Range("K12").Value = 61.11501429253024
$ws.Range("C13").Value = 7056.523220242473
$ws.Range("F13").Value = -2297.963454271704
$ws.Range("G13").Value = 15077.63455762254
$ws.Range("K13").Value = 48.79778864770837
$ws.Range("C14").Value = 6771.468230732743
$ws.Range("F14").Value = -1987.510904179602
$ws.Range("G14").Value = 11678.79904100106
$ws.Range("K14").Value = 37.82721141996331
$ws.Range("C15").Value = 6077.207375224643
$ws.Range("F15").Value = -1701.676930506009
$ws.Range("G15").Value = 9012.017504123072
$ws.Range("K15").Value = 28.4975192406987
$ws.Range("C16").Value = 5233.909756217452
$ws.Range("F16").Value = -1426.138684620535
$ws.Range("G16").Value = 6616.189777221574
$ws.Range("K16").Value = 20.49171240990674
$ws.Range("C17").Value = 4320.12641280349
$ws.Range("F17").Value = -1157.204843849326
$ws.Range("G17").Value = 4608.298301604259
$ws.Range("K17").Value = 13.87146478033595
$ws.Range("C18").Value = 3336.680672534461
$ws.Range("F18").Value = -897.9716041544795
$ws.Range("G18").Value = 2979.044593351941
$ws.Range("K18").Value = 8.596543661801414
$ws.Range("C19").Value = 2273.097302460856
$ws.Range("F19").Value = -648.4794598771598
$ws.Range("G19").Value = 1715.027634685314
$ws.Range("K19").Value = 4.621080887771636
$ws.Range("C20").Value = 1178.590700430782
$ws.Range("F20").Value = -397.2156585554874
$ws.Range("G20").Value = 803.1070853603685
$ws.Range("K20").Value = 1.929094280708807
$ws.Range("C21").Value = 337.5814046661731
$ws.Range("F21").Value = -175.0694776766458
$ws.Range("G21").Value = 245.6809126370629
$ws.Range("K21").Value = 0.4529288330461948
